# Translate row 1 header labels from Dutch to English
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value  = "Sample description"
$ws.Range("E1").Value  = "Iron II"
$ws.Range("F1").Value  = "MN II"
$ws.Range("G1").Value  = "cyanid"
$ws.Range("I1").Value  = "phosphor"
$ws.Range("J1").Value  = "benzene"
$ws.Range("K1").Value  = "toluene"
$ws.Range("L1").Value  = "ethylbenzene"
$ws.Range("M1").Value  = "o-xylene"
$ws.Range("N1").Value  = "(m+p)-xylene"
$ws.Range("O1").Value  = "sum xylenes (factor 0.7)"
$ws.Range("P1").Value  = "total BTEX (factor 0.7)"
$ws.Range("Q1").Value  = "naphthalene"
$ws.Range("R1").Value  = "naphthalene"
$ws.Range("S1").Value  = "acenaphthylene"
$ws.Range("T1").Value  = "acenaphtene"
$ws.Range("U1").Value  = "fluorene"
$ws.Range("V1").Value  = "phenanthrene"
$ws.Range("W1").Value  = "anthracene"
$ws.Range("X1").Value  = "fluoranthene"
$ws.Range("Y1").Value  = "pyrene"
$ws.Range("AA1").Value = "chrysene"
$ws.Range("AB1").Value = "benzo(b)fluoranthene"
$ws.Range("AC1").Value = "benzo(k)fluoranthene"
$ws.Range("AD1").Value = "benzo(a)pyrene"
$ws.Range("AE1").Value = "dibenz(a,h)anthracene"
$ws.Range("AF1").Value = "benzo(g,h,i)perylene"
$ws.Range("AG1").Value = "indeno(1,2,3-cd)pyrene"
$ws.Range("AH1").Value = "sum PAH (16 EPA)"
$ws.Range("AI1").Value = "sum PAH (VROM) (factor 0.7)"
$ws.Range("AJ1").Value = "fraction C10-C12"
$ws.Range("AK1").Value = "fraction C12-C22"
$ws.Range("AL1").Value = "fraction C22-C30"
$ws.Range("AM1").Value = "fraction C30-C40"
$ws.Range("AN1").Value = "total oil C10 - C40"
$ws.Range("AP1").Value = "nitrite"
$ws.Range("AQ1").Value = "nitrite"
$ws.Range("AR1").Value = "nitrate"
$ws.Range("AS1").Value = "nitrate"
$ws.Range("AT1").Value = "sulphates"
$ws.Range("AU1").Value = "Oxygen"
